# "add: prediksi ternak API" — refresh the sample production-data sheet:
#   - Ternak ID 1's production date moves from 1/1/2021 to 7/6/2023
#   - "Melahirkan" (calving) status relabeled to "Laktasi" (lactation)
#     across all rows (columns I and M)
#   - "Data Harian Rata Rata" (column E) values bumped up with new
#     highlight colors for each row

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Ternak ID 1's production date: keep it as literal text (not an
# auto-converted date serial) by leading with an apostrophe, just like
# the existing "d/m/yyyy"-style text already in that column.
$ws.Range("B2").Value = "'7/6/2023"

# Column E: "Data Harian Rata Rata" values + new fill colors per row.
$ws.Range("E2").Value = 1800
$ws.Range("E2").Interior.Color = 16774144   # RGB(0x00,0xF4,0xFF)

$ws.Range("E3").Value = 3050
$ws.Range("E3").Interior.Color = 16748544   # RGB(0x00,0x90,0xFF)

$ws.Range("E4").Value = 2200
$ws.Range("E4").Interior.Color = 16762880   # RGB(0x00,0xC8,0xFF)

$ws.Range("E5").Value = 950
$ws.Range("E5").Interior.Color = 36351      # RGB(0xFF,0x8D,0x00)

# Columns I & M ("Melahirkan" -> "Laktasi") for every data row, with the
# highlight fill switched from green to white.
$labelRanges = @("I2", "M2", "I3", "M3", "I4", "M4", "I5", "M5")
foreach ($addr in $labelRanges) {
    $cell = $ws.Range($addr)
    $cell.Value = "Laktasi"
    $cell.Interior.Color = 16777215   # RGB(0xFF,0xFF,0xFF)
}
